$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-02 Tuesday", "2025-09-03 Wednesday"),
    @("51×84=", "33×47="),
    @("51×45=", "13×39="),
    @("53×90=", "42×50="),
    @("57×11=", "87×51="),
    @("81×84=", "59×79="),
    @("26×97=", "23×29="),
    @("96×75=", "34×24="),
    @("24×57=", "13×42="),
    @("84×48=", "55×65="),
    @("21×35=", "92×11="),
    @("54×15=", "47×14="),
    @("39×49=", "83×55="),
    @("20×65=", "19×28="),
    @("91×97=", "16×84="),
    @("30×37=", "95×57="),
    @("89×80=", "95×35="),
    @("50×12=", "14×57="),
    @("54×96=", "93×22="),
    @("35×51=", "14×96="),
    @("15×87=", "39×37="),
    @("17×45=", "34×12="),
    @("53×47=", "23×40="),
    @("45×97=", "67×30="),
    @("27×27=", "28×22="),
    @("87×42=", "65×49=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
